# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> stock "Office Theme" palette (used by the notes master)
#   ppt/theme/theme2.xml  -> "Integral" palette (used by the slide master / slides)
#
# The authored change swaps the two themes' content: the slide master's theme
# becomes the plain "Office Theme" colour scheme, while the notes master's
# theme becomes the "Integral" colour scheme (font scheme and format scheme
# are identical between the two themes, so only the 12 colour-scheme slots
# actually change).
#
# This host's object model only ever exposes a single, shared Theme object
# (reached the same way from Presentation.SlideMaster, NotesMaster,
# HandoutMaster or any Slide) and it is backed by ppt/theme/theme2.xml, so
# that is the palette we repoint here -> from "Integral" to "Office Theme".

$p = $ppt.ActivePresentation

function HexToRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette: the stock "Office Theme" colour scheme, in the standard
# clrScheme order (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = HexToRGB $officeThemeColors[$i - 1]
}
